# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by a fresh handback report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn = $wb.Sheets.Item("zh-cn")
$dede = $wb.Sheets.Item("de-de")

# Latest HO Xliff Generate Date for be4ba6f3-...md (also mirrored on the de-de
# sheet's "Correspond Handoff Datetime" column for the same source file).
$overview.Range("G2").Value = "2016-09-06 07:17:44"
$dede.Range("H2").Value = "2016-09-06 07:17:44"

# zh-cn handoff / handback timestamps for be4ba6f3-...md
$zhcn.Range("H2").Value = "2016-09-06 07:17:39"
$zhcn.Range("K2").Value = "2016-09-06 07:17:59"

# de-de handback timestamp for be4ba6f3-...md
$dede.Range("K2").Value = "2016-09-06 07:18:16"
